# Update "Home win" sheet: remove rows 3-7, refresh row 2 with latest prediction
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Home win")
$ws.Range("A3:F7").EntireRow.Delete()
$ws.Range("A2").Value = "03-01-2025 17:30"
$ws.Range("B2").Value = "CYPRUS"
$ws.Range("C2").Value = "1. DIVISION"
$ws.Range("D2").Value = "AEK Larnaca - Omonia Nicosia"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.95

# Update "Draw" sheet: refresh rows 2-4 with latest predictions
$ws = $wb.Worksheets.Item("Draw")
$ws.Range("A2").Value = "03-01-2025 14:00"
$ws.Range("B2").Value = "ALGERIA"
$ws.Range("C2").Value = "COUPE NATIONALE"
$ws.Range("D2").Value = "MO Bejaia - Khroub"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 3.1

$ws.Range("A3").Value = "03-01-2025 23:00"
$ws.Range("B3").Value = "TRINIDAD-AND-TOBAGO"
$ws.Range("C3").Value = "PRO LEAGUE"
$ws.Range("D3").Value = "Prison Service - San Juan Jabloteh"
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 3.6

$ws.Range("A4").Value = "03-01-2025 19:00"
$ws.Range("B4").Value = "WORLD"
$ws.Range("C4").Value = "CAF CHAMPIONS LEAGUE"
$ws.Range("D4").Value = "CR Belouizdad - Al Ahly"
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 2.95

# Update "Btts" sheet: remove rows 5-6, refresh rows 2-4
$ws = $wb.Worksheets.Item("Btts")
$ws.Range("A5:F6").EntireRow.Delete()

$ws.Range("A2").Value = "03-01-2025 20:00"
$ws.Range("B2").Value = "FRANCE"
$ws.Range("C2").Value = "LIGUE 1"
$ws.Range("D2").Value = "Nice - Rennes"
$ws.Range("E2").Value = 83.3
$ws.Range("F2").Value = 1.75

$ws.Range("A3").Value = "03-01-2025 19:00"
$ws.Range("B3").Value = "FRANCE"
$ws.Range("C3").Value = "LIGUE 2"
$ws.Range("D3").Value = "Rodez - RED Star FC 93"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.7

$ws.Range("A4").Value = "03-01-2025 18:00"
$ws.Range("B4").Value = "SPAIN"
$ws.Range("C4").Value = "COPA DEL REY"
$ws.Range("D4").Value = "Racing Ferrol - Rayo Vallecano"
$ws.Range("E4").Value = 76
$ws.Range("F4").Value = 2.2

# Update "Over_Under" sheet: remove rows 6-7, refresh rows 2-5
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Range("A6:H7").EntireRow.Delete()

$ws.Range("A2").Value = "03-01-2025 19:00"
$ws.Range("B2").Value = "FRANCE"
$ws.Range("C2").Value = "LIGUE 2"
$ws.Range("D2").Value = "Rodez - RED Star FC 93"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.7
$ws.Range("G2").Value = 45
$ws.Range("H2").Value = 2.7

$ws.Range("A3").Value = "03-01-2025 19:00"
$ws.Range("B3").Value = "ITALY"
$ws.Range("C3").Value = "SUPER CUP"
$ws.Range("D3").Value = "Juventus - AC Milan"
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 2.1
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 4

$ws.Range("A4").Value = "03-01-2025 17:00"
$ws.Range("B4").Value = "TURKEY"
$ws.Range("C4").Value = "SÜPER LIG"
$ws.Range("D4").Value = "Rizespor - Besiktas"
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 1.8
$ws.Range("G4").Value = 50

$ws.Range("A5").Value = "03-01-2025 17:00"
$ws.Range("B5").Value = "TURKEY"
$ws.Range("C5").Value = "SÜPER LIG"
$ws.Range("D5").Value = "Alanyaspor - Konyaspor"
$ws.Range("E5").Value = 80
$ws.Range("F5").Value = 2.1
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = 3.75

# Update "Away Win" sheet: refresh row 2 with latest prediction
$ws = $wb.Worksheets.Item("Away Win")
$ws.Range("A2").Value = "03-01-2025 10:45"
$ws.Range("B2").Value = "AUSTRALIA"
$ws.Range("C2").Value = "A-LEAGUE"
$ws.Range("D2").Value = "Perth Glory - Western United"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.83
